$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for case with 380 kV (rows 2-25, columns B-F and I-N)
$data = @{
    2 = @{ "B" = 1.02; "C" = 1.037444437602331; "D" = 1.042138589983732; "E" = 1.041094892765427; "F" = 1.050016610925958; "I" = 1.034075243823482; "J" = 1.042547274773737; "K" = 1.044916123799758; "L" = 1.043875379655358; "M" = 1.052772059767402; "N" = 1.017985379744404 }
    3 = @{ "B" = 1.02; "C" = 1.038802706682118; "D" = 1.043162929669785; "E" = 1.04240240476731; "F" = 1.051262755441174; "I" = 1.034328747927642; "J" = 1.043547872266179; "K" = 1.045750674159067; "L" = 1.044992144075422; "M" = 1.053829447732489; "N" = 1.01832817739493 }
    4 = @{ "B" = 1.02; "C" = 1.039680865649978; "D" = 1.043824911905987; "E" = 1.043248093144826; "F" = 1.052068655914; "I" = 1.034491047235553; "J" = 1.044194160872374; "K" = 1.046289247941697; "L" = 1.045713865296643; "M" = 1.054512652465019; "N" = 1.018549321356315 }
    5 = @{ "B" = 1.02; "C" = 1.040049872859116; "D" = 1.044103011894805; "E" = 1.043603538091141; "F" = 1.052407354581509; "I" = 1.034558863316713; "J" = 1.044465584633633; "K" = 1.046515322788206; "L" = 1.046017064601658; "M" = 1.054799636038349; "N" = 1.018642130974632 }
    6 = @{ "B" = 1.02; "C" = 1.040111820904117; "D" = 1.044149694601701; "E" = 1.04366321413052; "F" = 1.052464217691636; "I" = 1.034570225648034; "J" = 1.044511141778746; "K" = 1.046553261784508; "L" = 1.046067960741949; "M" = 1.0548478081013; "N" = 1.018657704801411 }
    7 = @{ "B" = 1.02; "C" = 1.039685797008096; "D" = 1.043828628664196; "E" = 1.043252842937119; "F" = 1.052073182018749; "I" = 1.034491955025295; "J" = 1.044197788730726; "K" = 1.046292270104244; "L" = 1.04571791749084; "M" = 1.054516488076298; "N" = 1.018550562107856 }
    8 = @{ "B" = 1.02; "C" = 1.037903625211054; "D" = 1.042484944473806; "E" = 1.041536849122719; "F" = 1.050437843904845; "I" = 1.034161276386954; "J" = 1.042885674329784; "K" = 1.045198462704642; "L" = 1.044252983691998; "M" = 1.053129616862515; "N" = 1.018101368867993 }
    9 = @{ "B" = 1.02; "C" = 1.034757407010874; "D" = 1.04011070858018; "E" = 1.038510136208234; "F" = 1.047552674579913; "I" = 1.03356525719108; "J" = 1.040564508828207; "K" = 1.043259934751922; "L" = 1.041664547282897; "M" = 1.05067802094264; "N" = 1.017304669807149 }
    10 = @{ "B" = 1.02; "C" = 1.032655747466325; "D" = 1.038523369996079; "E" = 1.036490142534406; "F" = 1.045626684206563; "I" = 1.033158904863437; "J" = 1.039010794082197; "K" = 1.041959977027243; "L" = 1.039933990800872; "M" = 1.049038236727926; "N" = 1.016770009925879 }
    11 = @{ "B" = 1.02; "C" = 1.031744651494427; "D" = 1.037834930011938; "E" = 1.035614893545387; "F" = 1.044792061599228; "I" = 1.032980802184243; "J" = 1.038336489250501; "K" = 1.041395244478994; "L" = 1.039183425341821; "M" = 1.048326874637706; "N" = 1.016537646969981 }
    12 = @{ "B" = 1.02; "C" = 1.031406064628091; "D" = 1.037579042666136; "E" = 1.03528969570684; "F" = 1.044481943141287; "I" = 1.032914322936971; "J" = 1.038085787943001; "K" = 1.041185197989271; "L" = 1.038904444143638; "M" = 1.048062440544322; "N" = 1.016451207924687 }
    13 = @{ "B" = 1.02; "C" = 1.03147870022665; "D" = 1.037633939115966; "E" = 1.03535945595533; "F" = 1.044548469302542; "I" = 1.03292859763905; "J" = 1.038139574894192; "K" = 1.041230266376988; "L" = 1.038964295119251; "M" = 1.048119171777062; "N" = 1.016469755253646 }
    14 = @{ "B" = 1.02; "C" = 1.031716667201365; "D" = 1.037813781771365; "E" = 1.035588014481102; "F" = 1.044766429223512; "I" = 1.032975313604183; "J" = 1.03831577100674; "K" = 1.041377887692764; "L" = 1.039160368522654; "M" = 1.048305020585231; "N" = 1.016530504531703 }
    15 = @{ "B" = 1.02; "C" = 1.031863264446157; "D" = 1.037924566138207; "E" = 1.035728824686759; "F" = 1.044900707830583; "I" = 1.033004053889701; "J" = 1.03842430006667; "K" = 1.041468804910723; "L" = 1.039281150816222; "M" = 1.048419501225647; "N" = 1.016567917025993 }
    16 = @{ "B" = 1.02; "C" = 1.032716191031623; "D" = 1.038569035805681; "E" = 1.036548217357902; "F" = 1.045682061205907; "I" = 1.033170679575745; "J" = 1.039055512795011; "K" = 1.041997417405234; "L" = 1.039983777247482; "M" = 1.049085419277279; "N" = 1.016785413048615 }
    17 = @{ "B" = 1.02; "C" = 1.033250920753904; "D" = 1.038972994779201; "E" = 1.037062042712313; "F" = 1.046172005145864; "I" = 1.033274623259062; "J" = 1.039451042009282; "K" = 1.042328506789307; "L" = 1.040424185964153; "M" = 1.049502774939665; "N" = 1.016921613731514 }
    18 = @{ "B" = 1.02; "C" = 1.033562717311781; "D" = 1.039208510138466; "E" = 1.03736169284721; "F" = 1.046457718121249; "I" = 1.033335044533144; "J" = 1.039681599548506; "K" = 1.042521448070027; "L" = 1.04068095118093; "M" = 1.04974608386137; "N" = 1.01700097518667 }
    19 = @{ "B" = 1.02; "C" = 1.033669014703214; "D" = 1.03928879671885; "E" = 1.037463856520826; "F" = 1.046555128287441; "I" = 1.033355611493058; "J" = 1.039760188725639; "K" = 1.042587206031195; "L" = 1.04076848163498; "M" = 1.049829024382767; "N" = 1.017028021500799 }
    20 = @{ "B" = 1.02; "C" = 1.033193559931879; "D" = 1.038929664902748; "E" = 1.037006919889246; "F" = 1.046119445347443; "I" = 1.03326349253746; "J" = 1.039408620799865; "K" = 1.042293002419144; "L" = 1.040376946485186; "M" = 1.049458009849236; "N" = 1.016907009187752 }
    21 = @{ "B" = 1.02; "C" = 1.031646596464455; "D" = 1.03776082731139; "E" = 1.035520712251078; "F" = 1.044702248319081; "I" = 1.032961565868808; "J" = 1.03826389214541; "K" = 1.041334424638426; "L" = 1.039102635008012; "M" = 1.048250298362834; "N" = 1.016512618956367 }
    22 = @{ "B" = 1.02; "C" = 1.030673000283526; "D" = 1.037024947162503; "E" = 1.034585743501989; "F" = 1.043810605551743; "I" = 1.032769857600627; "J" = 1.037542797395142; "K" = 1.04073010917172; "L" = 1.038300336434515; "M" = 1.047489788627486; "N" = 1.016263902675927 }
    23 = @{ "B" = 1.02; "C" = 1.031189214762094; "D" = 1.037415145590101; "E" = 1.035081439967433; "F" = 1.044283339933689; "I" = 1.032871663861415; "J" = 1.037925193364106; "K" = 1.04105062272233; "L" = 1.038725754598577; "M" = 1.047893061649802; "N" = 1.016395823117688 }
    24 = @{ "B" = 1.02; "C" = 1.033219479133876; "D" = 1.038949244139548; "E" = 1.037031827691528; "F" = 1.046143195052551; "I" = 1.03326852267198; "J" = 1.039427789575387; "K" = 1.042309045864245; "L" = 1.040398292335919; "M" = 1.049478237660744; "N" = 1.016913608606065 }
    25 = @{ "B" = 1.02; "C" = 1.035571496212434; "D" = 1.04072529019292; "E" = 1.039292982478276; "F" = 1.048298994098601; "I" = 1.033720926168759; "J" = 1.041165677454444; "K" = 1.0437624207158; "L" = 1.042334574622193; "M" = 1.051312753829433; "N" = 1.017511253186763 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $data[$row][$col]
    }
}

Write-Host "Updated vm_pu values for case with 380 kV"
